$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Range("B2").Value = 'activityDTOModelMapper'
$ws.Range("D2").Value = 'org.andante.activity.controller.mapper.ActivityDTOModelMapper'
$ws.Range("B4").Value = 'operationHttpStatusMapper'
$ws.Range("D4").Value = 'org.andante.mappers.OperationHttpStatusMapper'
$ws.Range("B6").Value = 'NULL_PAGE_SIZE_ERROR_MESSAGE'
$ws.Range("B7").Value = 'ACTIVITY_EMAIL_BLANK_MESSAGE'
$ws.Range("B8").Value = 'IDENTIFIERS_LIST_NULL_MESSAGE'
$ws.Range("B9").Value = 'NEGATIVE_PAGE_ERROR_MESSAGE'
$ws.Range("B10").Value = 'NULL_PAGE_ERROR_MESSAGE'
$ws.Range("D10").Value = 'java.lang.String'
$ws.Range("B12").Value = 'ACTIVITY_IDENTIFIER_NOT_BLANK_MESSAGE'
$ws.Range("B19").Value = 'eventTimestamp'
$ws.Range("D19").Value = 'java.time.LocalDateTime'
$ws.Range("B20").Value = 'affectsAll'
$ws.Range("D20").Value = 'java.lang.Boolean'
$ws.Range("B21").Value = 'priority'
$ws.Range("D21").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B22").Value = 'acknowledgedUsers'
$ws.Range("D22").Value = 'java.util.Set'
$ws.Range("B23").Value = 'affectedUsers'
$ws.Range("B24").Value = 'domain'
$ws.Range("D24").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B25").Value = 'relatedId'
$ws.Range("D25").Value = 'java.lang.String'
$ws.Range("B26").Value = 'id'
$ws.Range("B31").Value = 'imageUrl'
$ws.Range("B32").Value = 'observingUsers'
$ws.Range("D32").Value = 'java.util.Set'
$ws.Range("B33").Value = 'key'
$ws.Range("B34").Value = 'communityImageUrl'
$ws.Range("D34").Value = 'java.lang.String'
$ws.Range("B36").Value = 'username'
$ws.Range("B43").Value = 'enableStartTLS'
$ws.Range("D43").Value = 'java.lang.String'
$ws.Range("B44").Value = 'host'
$ws.Range("B45").Value = 'username'
$ws.Range("B46").Value = 'password'
$ws.Range("B47").Value = 'smtpAuth'
$ws.Range("B48").Value = 'port'
$ws.Range("D48").Value = 'java.lang.Integer'
$ws.Range("B49").Value = 'jwtTokenDecoder'
$ws.Range("D49").Value = 'org.andante.activity.controller.decoder.JWTTokenDecoder'
$ws.Range("B51").Value = 'userProfileService'
$ws.Range("D51").Value = 'org.andante.activity.logic.UserProfileService'
$ws.Range("B52").Value = 'IDENTIFIERS_LIST_NULL_ERROR_MESSAGE'
$ws.Range("B53").Value = 'IMAGE_URL_BLANK_ERROR_MESSAGE'
$ws.Range("B54").Value = 'profileService'
$ws.Range("D54").Value = 'org.andante.activity.logic.ProfileService'
$ws.Range("B55").Value = 'USERNAME_BLANK_ERROR_MESSAGE'
$ws.Range("D55").Value = 'java.lang.String'
$ws.Range("B56").Value = 'IDENTIFIERS_LIST_SIZE_ERROR_MESSAGE'
$ws.Range("D56").Value = 'java.lang.String'
$ws.Range("B57").Value = 'userProfileMapper'
$ws.Range("D57").Value = 'org.andante.activity.controller.mapper.UserProfileDTOModelMapper'
$ws.Range("B59").Value = 'key'
$ws.Range("B61").Value = 'communityImageUrl'
$ws.Range("B62").Value = 'observedUsers'
$ws.Range("D62").Value = 'java.util.Set'
$ws.Range("B63").Value = 'imageUrl'
$ws.Range("D63").Value = 'java.lang.String'
$ws.Range("B64").Value = 'username'
$ws.Range("B66").Value = 'affectedUsers'
$ws.Range("D66").Value = 'java.util.Set'
$ws.Range("B67").Value = 'domain'
$ws.Range("D67").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B68").Value = 'acknowledgedUsers'
$ws.Range("D68").Value = 'java.util.Set'
$ws.Range("B69").Value = 'priority'
$ws.Range("D69").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B70").Value = 'affectsAll$value'
$ws.Range("D70").Value = 'java.lang.Boolean'
$ws.Range("B71").Value = 'affectsAll$set'
$ws.Range("D71").Value = 'boolean'
$ws.Range("B72").Value = 'eventTimestamp'
$ws.Range("D72").Value = 'java.time.LocalDateTime'
$ws.Range("B73").Value = 'description'
$ws.Range("D73").Value = 'java.lang.String'
$ws.Range("B75").Value = 'key'
$ws.Range("D75").Value = 'java.lang.String'
$ws.Range("B76").Value = 'emailAddress'
$ws.Range("D76").Value = 'java.lang.String'
$ws.Range("B77").Value = 'subscriptionDate'
$ws.Range("D77").Value = 'java.time.LocalDateTime'
$ws.Range("B78").Value = 'isConfirmed'
$ws.Range("D78").Value = 'java.lang.Boolean'
$ws.Range("B81").Value = 'adminUsername'
$ws.Range("D81").Value = 'java.lang.String'
$ws.Range("B84").Value = 'adminPassword'
$ws.Range("B85").Value = 'keycloakAdminTokenPath'
$ws.Range("B86").Value = 'userProfileService'
$ws.Range("D86").Value = 'org.andante.activity.logic.UserProfileService'
$ws.Range("B87").Value = 'keycloakGetUserPath'
$ws.Range("B88").Value = 'emailAddress'
$ws.Range("D88").Value = 'java.lang.String'
$ws.Range("B89").Value = 'subscriptionDate'
$ws.Range("D89").Value = 'java.time.LocalDateTime'
$ws.Range("B90").Value = 'isConfirmed'
$ws.Range("D90").Value = 'java.lang.Boolean'
$ws.Range("B91").Value = 'emailAddress'
$ws.Range("D91").Value = 'java.lang.String'
$ws.Range("B92").Value = 'subscriptionDate'
$ws.Range("D92").Value = 'java.time.LocalDateTime'
$ws.Range("B96").Value = 'observed'
$ws.Range("D96").Value = 'java.util.Set'
$ws.Range("B97").Value = 'communityImageUrl'
$ws.Range("D97").Value = 'java.lang.String'
$ws.Range("B98").Value = 'id'
$ws.Range("D98").Value = 'java.lang.String'
$ws.Range("B99").Value = 'observers'
$ws.Range("D99").Value = 'java.util.Set'
$ws.Range("B103").Value = 'NEWSLETTER_TEMPLATE'
$ws.Range("B104").Value = 'NEWSLETTER_TITLE'
$ws.Range("D104").Value = 'java.lang.String'
$ws.Range("B105").Value = 'logo'
$ws.Range("D105").Value = 'org.springframework.core.io.Resource'
$ws.Range("B106").Value = 'templateEngine'
$ws.Range("D106").Value = 'org.thymeleaf.TemplateEngine'
$ws.Range("B107").Value = 'sender'
$ws.Range("B108").Value = 'mailSender'
$ws.Range("D108").Value = 'org.springframework.mail.javamail.JavaMailSender'
$ws.Range("B110").Value = 'userProfileRepository'
$ws.Range("D110").Value = 'org.andante.activity.repository.UserProfileRepository'
$ws.Range("B111").Value = 'USER_NOT_FOUND_EXCEPTION_MESSAGE'
$ws.Range("D111").Value = 'java.lang.String'
$ws.Range("B112").Value = 'USER_CONFLICT_EXCEPTION_MESSAGE'
$ws.Range("D112").Value = 'java.lang.String'
$ws.Range("B113").Value = 'userProfileModelEntityMapper'
$ws.Range("D113").Value = 'org.andante.activity.logic.mapper.UserProfileModelEntityMapper'
$ws.Range("B115").Value = 'NEWSLETTER_NOT_FOUND_EXCEPTION_MESSAGE'
$ws.Range("B116").Value = 'NEWSLETTER_CONFLICT_EXCEPTION_MESSAGE'
$ws.Range("D116").Value = 'java.lang.String'
$ws.Range("B117").Value = 'newsletterMapper'
$ws.Range("D117").Value = 'org.andante.activity.logic.mapper.NewsletterModelEntityMapper'
$ws.Range("B118").Value = 'newsletterRepository'
$ws.Range("D118").Value = 'org.andante.activity.repository.NewsletterRepository'
$ws.Range("B119").Value = 'id'
$ws.Range("B120").Value = 'observed'
$ws.Range("D120").Value = 'java.util.Set'
$ws.Range("B121").Value = 'imageUrl'
$ws.Range("D121").Value = 'java.lang.String'
$ws.Range("B122").Value = 'username'
$ws.Range("D122").Value = 'java.lang.String'
$ws.Range("B123").Value = 'communityImageUrl'
$ws.Range("B124").Value = 'observers'
$ws.Range("D124").Value = 'java.util.Set'
$ws.Range("B129").Value = 'privateToken'
$ws.Range("B130").Value = 'databaseId'
$ws.Range("B132").Value = 'activityModelEntityMapper'
$ws.Range("D132").Value = 'org.andante.activity.logic.mapper.ActivityModelEntityMapper'
$ws.Range("B134").Value = 'ACTIVITY_NOT_FOUND_EXCEPTION_MESSAGE'
$ws.Range("D134").Value = 'java.lang.String'
$ws.Range("B135").Value = 'rsqlParser'
$ws.Range("D135").Value = 'cz.jirutka.rsql.parser.RSQLParser'
$ws.Range("B137").Value = 'activityRepository'
$ws.Range("D137").Value = 'org.andante.activity.repository.ActivityRepository'
$ws.Range("B138").Value = 'rsqlVisitor'
$ws.Range("D138").Value = 'org.andante.rsql.PersistentRSQLVisitor'
$ws.Range("B140").Value = 'acknowledgedUsers'
$ws.Range("D140").Value = 'java.util.Set'
$ws.Range("B141").Value = 'description'
$ws.Range("D141").Value = 'java.lang.String'
$ws.Range("B142").Value = 'priority'
$ws.Range("D142").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B144").Value = 'affectsAll'
$ws.Range("D144").Value = 'java.lang.Boolean'
$ws.Range("B145").Value = 'affectedUsers'
$ws.Range("D145").Value = 'java.util.Set'
$ws.Range("B146").Value = 'domain'
$ws.Range("D146").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B147").Value = 'eventTimestamp'
$ws.Range("D147").Value = 'java.time.LocalDateTime'
$ws.Range("B148").Value = 'id'
$ws.Range("D148").Value = 'java.lang.String'
$ws.Range("B151").Value = 'newsletterMapper'
$ws.Range("D151").Value = 'org.andante.activity.controller.mapper.NewsletterDTOModelMapper'
$ws.Range("B152").Value = 'emailSender'
$ws.Range("D152").Value = 'org.andante.activity.controller.email.EmailSender'
$ws.Range("B153").Value = 'EMAIL_BLANK_ERROR_MESSAGE'
$ws.Range("D153").Value = 'java.lang.String'
$ws.Range("B154").Value = 'newsletterService'
$ws.Range("D154").Value = 'org.andante.activity.logic.NewsletterService'
$ws.Range("B164").Value = 'acknowledgedUsers'
$ws.Range("D164").Value = 'java.util.Set'
$ws.Range("B165").Value = 'priority'
$ws.Range("D165").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B166").Value = 'affectsAll'
$ws.Range("D166").Value = 'java.lang.Boolean'
$ws.Range("B167").Value = 'key'
$ws.Range("D167").Value = 'java.lang.String'
$ws.Range("B168").Value = 'description'
$ws.Range("D168").Value = 'java.lang.String'
$ws.Range("B169").Value = 'domain'
$ws.Range("D169").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B170").Value = 'eventTimestamp'
$ws.Range("D170").Value = 'java.time.LocalDateTime'
$ws.Range("B171").Value = 'relatedId'
$ws.Range("D171").Value = 'java.lang.String'
$ws.Range("B172").Value = 'affectedUsers'
$ws.Range("D172").Value = 'java.util.Set'
$ws.Range("B173").Value = 'TREEIFY_THRESHOLD'
$ws.Range("D173").Value = 'int'
$ws.Range("B174").Value = 'values'
$ws.Range("D174").Value = 'java.util.Collection'
$ws.Range("B175").Value = 'MAXIMUM_CAPACITY'
$ws.Range("D175").Value = 'int'
$ws.Range("B176").Value = 'val$productOutputDTO'
$ws.Range("D176").Value = 'org.andante.product.dto.ProductOutputDTO'
$ws.Range("B177").Value = 'this$0'
$ws.Range("D177").Value = 'org.andante.activity.logic.impl.DefaultRecommendationService'
$ws.Range("B178").Value = 'DEFAULT_LOAD_FACTOR'
$ws.Range("D178").Value = 'float'
$ws.Range("B179").Value = 'modCount'
$ws.Range("D179").Value = 'int'
$ws.Range("B181").Value = 'keySet'
$ws.Range("D181").Value = 'java.util.Set'
$ws.Range("B182").Value = 'size'
$ws.Range("B183").Value = 'entrySet'
$ws.Range("D183").Value = 'java.util.Set'
$ws.Range("B185").Value = 'UNTREEIFY_THRESHOLD'
$ws.Range("B186").Value = 'table'
$ws.Range("D186").Value = 'java.util.HashMap$Node[]'
$ws.Range("B188").Value = 'DEFAULT_INITIAL_CAPACITY'
$ws.Range("D188").Value = 'int'
